# The underlying data export got re-sorted: several whole observation
# records (i.e. all columns A:AY of a row) traded places with each other
# while keeping the same row numbers in the sheet. Reproduce that by
# swapping/rotating the full row contents using COM Range.Value array
# read/write.
#
# Columns Y and AA hold a literal date string ("2026-01-24") stored as
# plain text; assigning that same text back through Range.Value would
# make Excel auto-coerce it into a real date serial number, which is not
# what the source file does. Since Y/AA are identical across every row
# touched by these swaps, we simply never write them - we copy the rest
# of each row in two chunks (A:X and AB:AY), skipping Y and AA, and
# include Z on its own since it's an ordinary time string (that one
# round-trips fine as text and does need to move with its row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowChunks($rowNum) {
    $h = @{}
    $h.Part1 = $ws.Range("A" + $rowNum + ":X" + $rowNum).Value()
    $h.PartZ = $ws.Range("Z" + $rowNum + ":Z" + $rowNum).Value()
    $h.Part2 = $ws.Range("AB" + $rowNum + ":AY" + $rowNum).Value()
    return $h
}

function Set-RowChunks($rowNum, $chunks) {
    $ws.Range("A" + $rowNum + ":X" + $rowNum).Value = $chunks.Part1
    $ws.Range("Z" + $rowNum + ":Z" + $rowNum).Value = $chunks.PartZ
    $ws.Range("AB" + $rowNum + ":AY" + $rowNum).Value = $chunks.Part2
}

function Swap-Rows($r1, $r2) {
    $v1 = Get-RowChunks $r1
    $v2 = Get-RowChunks $r2
    Set-RowChunks $r1 $v2
    Set-RowChunks $r2 $v1
}

# Simple pairwise swaps: each pair of rows fully exchanged their data.
Swap-Rows 3 4
Swap-Rows 6 7
Swap-Rows 13 14
Swap-Rows 16 17

# Rows 25-27 form a 3-way rotation instead of a simple swap:
#   new row25 <- old row26
#   new row26 <- old row27
#   new row27 <- old row25
$v25 = Get-RowChunks 25
$v26 = Get-RowChunks 26
$v27 = Get-RowChunks 27

Set-RowChunks 25 $v26
Set-RowChunks 26 $v27
Set-RowChunks 27 $v25
